$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (M2:T2)
$ws.Range("M2").Value = 48.42420966666666
$ws.Range("N2").Value = 145.272629
$ws.Range("O2").Value = 0.6311762527593259
$ws.Range("P2").Value = 0.6311762527593258
$ws.Range("Q2").Value = 18.30228515438755
$ws.Range("R2").Value = 164.720566389488
$ws.Range("S2").Value = 0.6311762527593259
$ws.Range("T2").Value = 0.6311762527593258

# Row 3 updates (M3:T3)
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("O3").Value = 0.08928392431779728
$ws.Range("P3").Value = 0.08928392431779726
$ws.Range("S3").Value = 0.08928392431779728
$ws.Range("T3").Value = 0.08928392431779726

# Row 4 updates (N4, O4:P4, Q4:R4, S4:T4)
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2795398229228769
$ws.Range("P4").Value = 0.2795398229228769
$ws.Range("Q4").Value = 8.105846075125331
$ws.Range("R4").Value = 72.95261467612799
$ws.Range("S4").Value = 0.2795398229228769
$ws.Range("T4").Value = 0.2795398229228769

$wb.Save()
